$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Lama5"
$ws.Range("C2").Value = "Bcam"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 19.63445233333333
$ws.Range("H2").Value = 58.903357
$ws.Range("I2").Value = 0.514089849859583
$ws.Range("J2").Value = 0.5140898498595828
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 28.01710266666667
$ws.Range("N2").Value = 84.051308
$ws.Range("O2").Value = 0.541609178372911
$ws.Range("P2").Value = 0.541609178372911
$ws.Range("Q2").Value = 550.100466826773
$ws.Range("R2").Value = 4950.904201440956
$ws.Range("S2").Value = 0.2784357811923019
$ws.Range("T2").Value = 0.2784357811923018

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Lama5"
$ws.Range("C3").Value = "Bcam"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 19.63445233333333
$ws.Range("H3").Value = 58.903357
$ws.Range("I3").Value = 0.514089849859583
$ws.Range("J3").Value = 0.5140898498595828
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.9744343333333333
$ws.Range("N3").Value = 2.923303
$ws.Range("O3").Value = 0.01883715760812509
$ws.Range("P3").Value = 0.01883715760812509
$ws.Range("Q3").Value = 19.13248446979678
$ws.Range("R3").Value = 172.192360228171
$ws.Range("S3").Value = 0.00968399152654233
$ws.Range("T3").Value = 0.009683991526542327

# Row 4: ECs -> M2
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Lama5"
$ws.Range("C4").Value = "Bcam"
$ws.Range("D4").Value = "M2"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 19.63445233333333
$ws.Range("H4").Value = 58.903357
$ws.Range("I4").Value = 0.514089849859583
$ws.Range("J4").Value = 0.5140898498595828
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.07953766666666666
$ws.Range("N4").Value = 0.238613
$ws.Range("O4").Value = 0.001537572632172427
$ws.Range("P4").Value = 0.001537572632172427
$ws.Range("Q4").Value = 1.561678524871222
$ws.Range("R4").Value = 14.055106723841
$ws.Range("S4").Value = 0.0007904504836217268
$ws.Range("T4").Value = 0.0007904504836217267

# Row 5: ECs -> sCs
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Lama5"
$ws.Range("C5").Value = "Bcam"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 19.63445233333333
$ws.Range("H5").Value = 58.903357
$ws.Range("I5").Value = 0.514089849859583
$ws.Range("J5").Value = 0.5140898498595828
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 22.65829733333333
$ws.Range("N5").Value = 67.974892
$ws.Range("O5").Value = 0.4380160913867915
$ws.Range("P5").Value = 0.4380160913867915
$ws.Range("Q5").Value = 444.8832589458271
$ws.Range("R5").Value = 4003.949330512444
$ws.Range("S5").Value = 0.225179626657117
$ws.Range("T5").Value = 0.2251796266571169

# Row 6: FAPs -> ECs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Lama5"
$ws.Range("C6").Value = "Bcam"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.736962
$ws.Range("H6").Value = 2.210886
$ws.Range("I6").Value = 0.01929591299519065
$ws.Range("J6").Value = 0.01929591299519064
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 28.01710266666667
$ws.Range("N6").Value = 84.051308
$ws.Range("O6").Value = 0.541609178372911
$ws.Range("P6").Value = 0.541609178372911
$ws.Range("Q6").Value = 20.647540015432
$ws.Range("R6").Value = 185.827860138888
$ws.Range("S6").Value = 0.01045084358328038
$ws.Range("T6").Value = 0.01045084358328038

# Row 7: FAPs -> FAPs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Lama5"
$ws.Range("C7").Value = "Bcam"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.736962
$ws.Range("H7").Value = 2.210886
$ws.Range("I7").Value = 0.01929591299519065
$ws.Range("J7").Value = 0.01929591299519064
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9744343333333333
$ws.Range("N7").Value = 2.923303
$ws.Range("O7").Value = 0.01883715760812509
$ws.Range("P7").Value = 0.01883715760812509
$ws.Range("Q7").Value = 0.718121075162
$ws.Range("R7").Value = 6.463089676457999
$ws.Range("S7").Value = 0.0003634801542830753
$ws.Range("T7").Value = 0.0003634801542830753

# Row 8: FAPs -> M2
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Lama5"
$ws.Range("C8").Value = "Bcam"
$ws.Range("D8").Value = "M2"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.736962
$ws.Range("H8").Value = 2.210886
$ws.Range("I8").Value = 0.01929591299519065
$ws.Range("J8").Value = 0.01929591299519064
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.07953766666666666
$ws.Range("N8").Value = 0.238613
$ws.Range("O8").Value = 0.001537572632172427
$ws.Range("P8").Value = 0.001537572632172427
$ws.Range("Q8").Value = 0.058616237902
$ws.Range("R8").Value = 0.5275461411179999
$ws.Range("S8").Value = 0.00002966886773418542
$ws.Range("T8").Value = 0.00002966886773418542

# Row 9: FAPs -> sCs
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Lama5"
$ws.Range("C9").Value = "Bcam"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.736962
$ws.Range("H9").Value = 2.210886
$ws.Range("I9").Value = 0.01929591299519065
$ws.Range("J9").Value = 0.01929591299519064
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 22.65829733333333
$ws.Range("N9").Value = 67.974892
$ws.Range("O9").Value = 0.4380160913867915
$ws.Range("P9").Value = 0.4380160913867915
$ws.Range("Q9").Value = 16.698304119368
$ws.Range("R9").Value = 150.284737074312
$ws.Range("S9").Value = 0.008451920389893003
$ws.Range("T9").Value = 0.008451920389893001

# Row 10: M2 -> ECs
$ws.Range("A10").Value = "M2"
$ws.Range("B10").Value = "Lama5"
$ws.Range("C10").Value = "Bcam"
$ws.Range("D10").Value = "ECs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1783006666666667
$ws.Range("H10").Value = 0.534902
$ws.Range("I10").Value = 0.004668455294824549
$ws.Range("J10").Value = 0.004668455294824548
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 28.01710266666667
$ws.Range("N10").Value = 84.051308
$ws.Range("O10").Value = 0.541609178372911
$ws.Range("P10").Value = 0.541609178372911
$ws.Range("Q10").Value = 4.995468083535112
$ws.Range("R10").Value = 44.959212751816
$ws.Range("S10").Value = 0.00252847823650059
$ws.Range("T10").Value = 0.002528478236500589

# Row 11: M2 -> FAPs
$ws.Range("A11").Value = "M2"
$ws.Range("B11").Value = "Lama5"
$ws.Range("C11").Value = "Bcam"
$ws.Range("D11").Value = "FAPs"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1783006666666667
$ws.Range("H11").Value = 0.534902
$ws.Range("I11").Value = 0.004668455294824549
$ws.Range("J11").Value = 0.004668455294824548
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 0.9744343333333333
$ws.Range("N11").Value = 2.923303
$ws.Range("O11").Value = 0.01883715760812509
$ws.Range("P11").Value = 0.01883715760812509
$ws.Range("Q11").Value = 0.1737422912562222
$ws.Range("R11").Value = 1.563680621306
$ws.Range("S11").Value = 0.00008794042817509613
$ws.Range("T11").Value = 0.0000879404281750961

# Row 12: M2 -> M2
$ws.Range("A12").Value = "M2"
$ws.Range("B12").Value = "Lama5"
$ws.Range("C12").Value = "Bcam"
$ws.Range("D12").Value = "M2"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1783006666666667
$ws.Range("H12").Value = 0.534902
$ws.Range("I12").Value = 0.004668455294824549
$ws.Range("J12").Value = 0.004668455294824548
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.07953766666666666
$ws.Range("N12").Value = 0.238613
$ws.Range("O12").Value = 0.001537572632172427
$ws.Range("P12").Value = 0.001537572632172427
$ws.Range("Q12").Value = 0.01418161899177778
$ws.Range("R12").Value = 0.127634570926
$ws.Range("S12").Value = 0.000007178089095842687
$ws.Range("T12").Value = 0.000007178089095842685

# Row 13: M2 -> sCs
$ws.Range("A13").Value = "M2"
$ws.Range("B13").Value = "Lama5"
$ws.Range("C13").Value = "Bcam"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1783006666666667
$ws.Range("H13").Value = 0.534902
$ws.Range("I13").Value = 0.004668455294824549
$ws.Range("J13").Value = 0.004668455294824548
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 22.65829733333333
$ws.Range("N13").Value = 67.974892
$ws.Range("O13").Value = 0.541609178372911
$ws.Range("P13").Value = 0.541609178372911
$ws.Range("Q13").Value = 4.039989520064889
$ws.Range("R13").Value = 36.359905680584
$ws.Range("S13").Value = 0.002044858541053021
$ws.Range("T13").Value = 0.00204485854105302

# Row 14: sCs -> ECs
$ws.Range("A14").Value = "sCs"
$ws.Range("B14").Value = "Lama5"
$ws.Range("C14").Value = "Bcam"
$ws.Range("D14").Value = "ECs"
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 17.64293233333333
$ws.Range("H14").Value = 52.928797
$ws.Range("I14").Value = 0.461945781850402
$ws.Range("J14").Value = 0.4619457818504019
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 28.01710266666667
$ws.Range("N14").Value = 84.051308
$ws.Range("O14").Value = 0.541609178372911
$ws.Range("P14").Value = 0.541609178372911
$ws.Range("Q14").Value = 494.303846524053
$ws.Range("R14").Value = 4448.734618716477
$ws.Range("S14").Value = 0.2501940753608282
$ws.Range("T14").Value = 0.2501940753608282

# Row 15: sCs -> FAPs
$ws.Range("A15").Value = "sCs"
$ws.Range("B15").Value = "Lama5"
$ws.Range("C15").Value = "Bcam"
$ws.Range("D15").Value = "FAPs"
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 17.64293233333333
$ws.Range("H15").Value = 52.928797
$ws.Range("I15").Value = 0.461945781850402
$ws.Range("J15").Value = 0.4619457818504019
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 0.9744343333333333
$ws.Range("N15").Value = 2.923303
$ws.Range("O15").Value = 0.01883715760812509
$ws.Range("P15").Value = 0.01883715760812509
$ws.Range("Q15").Value = 17.19187900627678
$ws.Range("R15").Value = 154.726911056491
$ws.Range("S15").Value = 0.008701745499124594
$ws.Range("T15").Value = 0.008701745499124591

# Row 16: sCs -> M2
$ws.Range("A16").Value = "sCs"
$ws.Range("B16").Value = "Lama5"
$ws.Range("C16").Value = "Bcam"
$ws.Range("D16").Value = "M2"
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 17.64293233333333
$ws.Range("H16").Value = 52.928797
$ws.Range("I16").Value = 0.461945781850402
$ws.Range("J16").Value = 0.4619457818504019
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.07953766666666666
$ws.Range("N16").Value = 0.238613
$ws.Range("O16").Value = 0.001537572632172427
$ws.Range("P16").Value = 0.001537572632172427
$ws.Range("Q16").Value = 1.403277670951222
$ws.Range("R16").Value = 12.629499038561
$ws.Range("S16").Value = 0.0007102751917206723
$ws.Range("T16").Value = 0.0007102751917206723

# Row 17: sCs -> sCs
$ws.Range("A17").Value = "sCs"
$ws.Range("B17").Value = "Lama5"
$ws.Range("C17").Value = "Bcam"
$ws.Range("D17").Value = "sCs"
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 17.64293233333333
$ws.Range("H17").Value = 52.928797
$ws.Range("I17").Value = 0.461945781850402
$ws.Range("J17").Value = 0.4619457818504019
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 22.65829733333333
$ws.Range("N17").Value = 67.974892
$ws.Range("O17").Value = 0.4380160913867915
$ws.Range("P17").Value = 0.4380160913867915
$ws.Range("Q17").Value = 399.7588066405472
$ws.Range("R17").Value = 3597.829259764924
$ws.Range("S17").Value = 0.2023396857987285
$ws.Range("T17").Value = 0.2023396857987285
